$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 97, shifting rows 97:151 down to 98:152
$ws.Rows("97:97").Insert()

# Populate the new row 97 with the new data
$ws.Range("A97").Value = 6
$ws.Range("B97").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C97").Value = "Metropolitana"
$ws.Range("D97").Value = 44582
$ws.Range("D97").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E97").Value = 13
$ws.Range("F97").Value = "Fruta"
$ws.Range("G97").Value = 100101
$ws.Range("H97").Value = "Berries"
$ws.Range("I97").Value = 100101004
$ws.Range("J97").Value = "Frambuesa"
$ws.Range("K97").Value = "Sin especificar"
$ws.Range("L97").Value = "Especial"
$ws.Range("M97").Value = 200
$ws.Range("N97").Value = 8000
$ws.Range("O97").Value = 8000
$ws.Range("P97").Value = 8000
$ws.Range("Q97").Value = '$/bandeja 2 kilos'
$ws.Range("R97").Value = "Provincia de Colchagua"
$ws.Range("S97").Value = 4000
$ws.Range("T97").Value = 2
